# The commit removes the "json database stuff" placeholder by appending the
# first real data row to the "nutrition" sheet (sheet2):
#   A2 = "213_8290"               (text)
#   B2 = "2022/01/23 13:41:26"    (text)
#   C2 = "2"                      (text, looks numeric but stored as text)
#   D2 = 3                        (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nutrition")

$ws.Range("A2").Value = "213_8290"
$ws.Range("B2").Value = "2022/01/23 13:41:26"

# C2 must stay text ("2"), not be auto-coerced to the number 2: force a text
# number format before assigning, then drop back to the sheet's default
# style so no stray cell style is left behind.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = 3
